# Refresh the crypto price/volume snapshot (GitHub Actions scheduled update).
# Price cells (column D) hold text that looks numeric (e.g. "299.47"), so we
# force the cell to Text format before writing the value - otherwise Excel's
# normal type inference would coerce the string to a real number and silently
# drop meaningful trailing zeros (e.g. "24.30" -> 24.3, "9.10" -> 9.1).
# Volume cells (column E) are padded with spaces ("  +0.73%  "), which already
# keeps Excel from treating them as numbers, so no special handling is needed.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "42.930.28"
$ws.Range("E2").Value = "  +0.73%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.297.14"
$ws.Range("E3").Value = "  +0.54%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "299.47"
$ws.Range("E5").Value = "  -0.47%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "97.10"
$ws.Range("E6").Value = "  -0.02%  "
$ws.Range("E7").Value = "  +0.65%  "
$ws.Range("E8").Value = "  -0.02%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.507"
$ws.Range("E9").Value = "  +1.16%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "33.74"
$ws.Range("E10").Value = "  +0.71%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0791"
$ws.Range("E11").Value = "  +0.51%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "48.99"
$ws.Range("E12").Value = "  -3.14%  "
$ws.Range("E13").Value = "  +2.89%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "17.06"
$ws.Range("E14").Value = "  +12.19%  "
$ws.Range("E15").Value = "  +1.89%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.647.88"
$ws.Range("E16").Value = "  +0.24%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.336.11"
$ws.Range("E17").Value = "  +1.63%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.807"
$ws.Range("E18").Value = "  +2.40%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "42.876.55"
$ws.Range("E19").Value = "  +0.76%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.63"
$ws.Range("E20").Value = "  +1.19%  "
$ws.Range("E21").Value = "  +0.70%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.04"
$ws.Range("E22").Value = "  +0.62%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "67.45"
$ws.Range("E23").Value = "  +1.11%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "236.51"
$ws.Range("E24").Value = "  +0.60%  "
$ws.Range("E25").Value = "  +5.10%  "
$ws.Range("E26").Value = "  +0.00%  "
$ws.Range("E27").Value = "  -1.48%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "24.30"
$ws.Range("E28").Value = "  -0.71%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "166.67"
$ws.Range("E29").Value = "  +0.51%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.07"
$ws.Range("E30").Value = "  +0.02%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "33.70"
$ws.Range("E31").Value = "  +0.12%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "9.11"
$ws.Range("E32").Value = "  +0.07%  "
$ws.Range("E33").Value = "  +0.00%  "
$ws.Range("E34").Value = "  -0.27%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.60"
$ws.Range("E35").Value = "  +5.41%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.44"
$ws.Range("E36").Value = "  +1.93%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "16.74"
$ws.Range("E37").Value = "  +3.06%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0699"
$ws.Range("E38").Value = "  +0.50%  "
$ws.Range("E39").Value = "  -0.26%  "
$ws.Range("E40").Value = "  +0.46%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.77"
$ws.Range("E41").Value = "  -0.02%  "
$ws.Range("E42").Value = "  -0.37%  "
$ws.Range("E43").Value = "  -1.88%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.988.25"
$ws.Range("E44").Value = "  +1.29%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0284"
$ws.Range("E45").Value = "  +0.67%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "9.84"
$ws.Range("E46").Value = "  +1.73%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "17.45"
$ws.Range("E47").Value = "  -2.06%  "
$ws.Range("E48").Value = "  +0.02%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.527.54"
$ws.Range("E49").Value = "  +0.82%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "53.12"
$ws.Range("E50").Value = "  -0.27%  "
$ws.Range("E51").Value = "  -2.07%  "
